$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D8").Value = 1
